$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue 'D2' '26.284.09'
Set-TextValue 'E2' '  +0.91%  '
Set-TextValue 'D3' '1.679.15'
Set-TextValue 'E3' '  +0.67%  '
Set-TextValue 'E4' '  +0.16%  '
Set-TextValue 'D5' '217.96'
Set-TextValue 'E5' '  +0.55%  '
Set-TextValue 'D6' '0.5341'
Set-TextValue 'E6' '  +4.50%  '
Set-TextValue 'E7' '  +0.15%  '
Set-TextValue 'D8' '0.2684'
Set-TextValue 'E8' '  +1.06%  '
Set-TextValue 'D9' '0.06471'
Set-TextValue 'E9' '  +1.01%  '
Set-TextValue 'D10' '21.95'
Set-TextValue 'E10' '  +0.46%  '
Set-TextValue 'D11' '0.07542'
Set-TextValue 'E11' '  +1.74%  '
Set-TextValue 'D12' '1.680.59'
Set-TextValue 'E12' '  -0.69%  '
Set-TextValue 'D13' '4.522'
Set-TextValue 'E13' '  +0.47%  '
Set-TextValue 'D14' '0.5779'
Set-TextValue 'E14' '  -0.99%  '
Set-TextValue 'D15' '0.000008468'
Set-TextValue 'E15' '  -0.93%  '
Set-TextValue 'D16' '64.76'
Set-TextValue 'E16' '  +0.64%  '
Set-TextValue 'D17' '26.317.02'
Set-TextValue 'E17' '  +0.85%  '
Set-TextValue 'D18' '4.902'
Set-TextValue 'E18' '  -0.91%  '
Set-TextValue 'D20' '10.87'
Set-TextValue 'E20' '  +0.96%  '
Set-TextValue 'D21' '191.29'
Set-TextValue 'E21' '  +0.57%  '
Set-TextValue 'D22' '6.204'
Set-TextValue 'E22' '  -0.41%  '
Set-TextValue 'D23' '1.007'
Set-TextValue 'E23' '  +0.10%  '
Set-TextValue 'D24' '145.89'
Set-TextValue 'E24' '  +0.45%  '
Set-TextValue 'D25' '7.827'
Set-TextValue 'E25' '  +2.57%  '
Set-TextValue 'D26' '0.1274'
Set-TextValue 'E26' '  +5.84%  '
Set-TextValue 'D27' '15.76'
Set-TextValue 'E27' '  +0.99%  '
Set-TextValue 'D28' '0.06499'
Set-TextValue 'E28' '  +0.17%  '
Set-TextValue 'D29' '1.378'
Set-TextValue 'E29' '  +4.16%  '
Set-TextValue 'D30' '1.320'
Set-TextValue 'E30' '  +0.26%  '
Set-TextValue 'D31' '3.582'
Set-TextValue 'E31' '  +1.12%  '
Set-TextValue 'E32' '  +1.89%  '
Set-TextValue 'D33' '1.664'
Set-TextValue 'E33' '  +0.85%  '
Set-TextValue 'D34' '1.032'
Set-TextValue 'E34' '  +1.40%  '
Set-TextValue 'D35' '0.6176'
Set-TextValue 'E35' '  +1.28%  '
Set-TextValue 'D37' '2.704'
Set-TextValue 'E37' '  -0.28%  '
Set-TextValue 'D38' '6.252'
Set-TextValue 'E38' '  +0.42%  '
Set-TextValue 'D39' '1.112.64'
Set-TextValue 'E39' '  +2.40%  '
Set-TextValue 'D40' '0.01624'
Set-TextValue 'E40' '  +1.28%  '
Set-TextValue 'D41' '0.8702'
Set-TextValue 'E41' '  +0.57%  '
Set-TextValue 'D42' '1.014'
Set-TextValue 'E42' '  +0.48%  '
Set-TextValue 'D43' '100.34'
Set-TextValue 'D44' '1.828.91'
Set-TextValue 'E44' '  +0.70%  '
Set-TextValue 'D45' '0.00000000110'
Set-TextValue 'E45' '  -4.65%  '
Set-TextValue 'D46' '57.17'
Set-TextValue 'E46' '  +1.56%  '
Set-TextValue 'D47' '8.179'
Set-TextValue 'E47' '  +1.26%  '
Set-TextValue 'D48' '1.004'
Set-TextValue 'E48' '  -0.50%  '
Set-TextValue 'D49' '0.05262'
Set-TextValue 'E49' '  +0.54%  '
Set-TextValue 'B50' 'Aptos'
Set-TextValue 'C50' 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 'D50' '6.083'
Set-TextValue 'E50' '  +1.04%  '
Set-TextValue 'B51' 'Mantle'
Set-TextValue 'C51' 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
Set-TextValue 'D51' '0.4288'
Set-TextValue 'E51' '  -0.01%  '
